$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text edits (rich-text runs) ---
$ws.Range("A8").Value = "Volume 30   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/29/2023  Through  6/4/2023"

# --- Convert numeric-styled cells to "N/A" / "***.*" text placeholders ---
# Use Copy(Destination) from a stable template cell so both the shared-string
# text AND the style id (14) transfer atomically.
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("M14").Copy($ws.Range("E14"))
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("M14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("M14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("M14").Copy($ws.Range("E29"))

# --- Convert "N/A" / "***.*" placeholder cells back into numeric cells ---
# First Copy(Destination) from a stable numeric-styled template cell to fix
# the style/type (15 = count, 16 = percent), then set the final numeric value.
$ws.Range("I15").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("H15").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("I15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$ws.Range("H15").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 100
$ws.Range("I15").Copy($ws.Range("F22"))
$ws.Range("F22").Value = 1
$ws.Range("I15").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("I15").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("H15").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 0
$ws.Range("I15").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("H15").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -25
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -19.047619047619
$ws.Range("I16").Value = 79
$ws.Range("J16").Value = 78
$ws.Range("K16").Value = 1.282051282051
$ws.Range("L16").Value = 14.492753623188
$ws.Range("M16").Value = -22.549019607843
$ws.Range("N16").Value = -58.421052631578
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -18.518518518518
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 15.094339622641
$ws.Range("L17").Value = 41.860465116279
$ws.Range("M17").Value = 58.441558441558
$ws.Range("N17").Value = 29.787234042553
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 63
$ws.Range("J18").Value = 42
$ws.Range("L18").Value = 53.658536585365
$ws.Range("M18").Value = -52.272727272727
$ws.Range("N18").Value = -83.2
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 76.470588235294
$ws.Range("I19").Value = 286
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = 27.111111111111
$ws.Range("L19").Value = 67.251461988304
$ws.Range("M19").Value = 70.238095238095
$ws.Range("N19").Value = 55.434782608695
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -37.5
$ws.Range("F20").Value = 48
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 152.631578947368
$ws.Range("I20").Value = 233
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = 113.761467889908
$ws.Range("L20").Value = 164.772727272727
$ws.Range("M20").Value = 219.178082191781
$ws.Range("N20").Value = -71.234567901234
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = 42.342342342342
$ws.Range("I21").Value = 790
$ws.Range("J21").Value = 570
$ws.Range("K21").Value = 38.596491228070
$ws.Range("L21").Value = 70.626349892008
$ws.Range("M21").Value = 41.071428571428
$ws.Range("N21").Value = -52.581032412965
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 8
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 14.285714285714
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 26
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 85.714285714285
$ws.Range("L23").Value = 44.444444444444
$ws.Range("M23").Value = 62.5
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 31.578947368421
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 14.285714285714
$ws.Range("I24").Value = 579
$ws.Range("J24").Value = 485
$ws.Range("K24").Value = 19.381443298969
$ws.Range("L24").Value = 51.570680628272
$ws.Range("M24").Value = 5.272727272727
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 34.482758620689
$ws.Range("I25").Value = 216
$ws.Range("J25").Value = 190
$ws.Range("K25").Value = 13.684210526315
$ws.Range("L25").Value = 24.137931034482
$ws.Range("M25").Value = 26.315789473684
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -27.777777777777
$ws.Range("L26").Value = 116.666666666667
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -60
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -29.166666666666
$ws.Range("L27").Value = -15
$ws.Range("N28").Value = -71.428571428571
$ws.Range("N29").Value = -71.428571428571
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = -33.333333333333
